# "added a case with 2 variants"
# Appends a new phishing-sample record (id 44) as row 45 of Sheet1,
# right below the existing last row (44), following the same column
# layout as the rest of the table:
#   A: id  B: type  C: added  D: source  E: lang  F: motivation
#   G: lang  H: personalised  I: description  J: Entity  (K: moreinfo - blank here)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 45

$ws.Cells.Item($newRow, 1).Value  = 44
$ws.Cells.Item($newRow, 2).Value  = "msg"
$ws.Cells.Item($newRow, 3).Value  = "2021-08-08"
$ws.Cells.Item($newRow, 4).Value  = "MCAST"
$ws.Cells.Item($newRow, 5).Value  = "shortened"
$ws.Cells.Item($newRow, 6).Value  = "delivery"
$ws.Cells.Item($newRow, 7).Value  = "mt"
$ws.Cells.Item($newRow, 8).Value  = "no"
$ws.Cells.Item($newRow, 9).Value  = "confirm address for postal delivery"
$ws.Cells.Item($newRow, 10).Value = "DHL"

# Matches the saved cursor position after entering the new row.
$ws.Range("J45").Select()
